# Add two new result sheets "t8_ML_1612" and "t9_ML_2012" at the end of the
# workbook, after the existing "t7_ML_0912" sheet, and populate them with the
# quiz results (same 14-column layout used by t6/t7: Nom de famille, Prénom,
# Clé, Adresse de courriel, Durée, Note/20, Q1..Q8).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data: header labels shared by both new sheets
# ---------------------------------------------------------------------
$headers = @(
    "Nom de famille",
    "Prénom",
    "Clé",
    "Adresse de courriel",
    "Durée",
    "Note/20,00",
    "Q. 1 /2,00",
    "Q. 2 /2,00",
    "Q. 3 /2,50",
    "Q. 4 /2,50",
    "Q. 5 /3,00",
    "Q. 6 /3,00",
    "Q. 7 /2,00",
    "Q. 8 /3,00"
)

# =======================================================================
# Sheet "t8_ML_1612" (3 students, numeric question scores)
# =======================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $lastSheet)
$ws8.Name = "t8_ML_1612"

for ($c = 1; $c -le $headers.Length; $c++) {
    $ws8.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$t8rows = @(
    @{ A="PECOURT";           B="Bertille";                       D="bertille.pecourt@etu.unilasalle.fr";            E="7 min 28 s";  F=4;   G=2; H=2; I=0; J=0;   K=0; L=0; M=0; N=0 },
    @{ A="CHIEN-CHOW-CHINE";  B="Jules";                          D="jules.chien-chow-chine@etu.unilasalle.fr";      E="10 min 8 s";  F=6;   G=0; H=0; I=0; J=0;   K=3; L=0; M=0; N=3 },
    @{ A="BOUTILLIER";        B="Hugo";                           D="hugo.boutillier@etu.unilasalle.fr";             E="7 min 33 s";  F=8.5; G=0; H=0; I=0; J=2.5; K=0; L=3; M=0; N=3 }
)

$r = 2
foreach ($row in $t8rows) {
    $ws8.Cells.Item($r, 1).Value = $row.A
    $ws8.Cells.Item($r, 2).Value = $row.B
    $ws8.Cells.Item($r, 3).Formula = "=A$r&B$r"
    $ws8.Cells.Item($r, 4).Value = $row.D
    $ws8.Cells.Item($r, 5).Value = $row.E
    $ws8.Cells.Item($r, 6).Value = $row.F
    $ws8.Cells.Item($r, 7).Value = $row.G
    $ws8.Cells.Item($r, 8).Value = $row.H
    $ws8.Cells.Item($r, 9).Value = $row.I
    $ws8.Cells.Item($r, 10).Value = $row.J
    $ws8.Cells.Item($r, 11).Value = $row.K
    $ws8.Cells.Item($r, 12).Value = $row.L
    $ws8.Cells.Item($r, 13).Value = $row.M
    $ws8.Cells.Item($r, 14).Value = $row.N
    $r++
}

[void]$ws8.Range("L2").Select()

# =======================================================================
# Sheet "t9_ML_2012" (6 students, question scores stored as French-format
# text, same convention as "t6_ML_0212" / "t7_ML_0912")
# =======================================================================
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add($null, $lastSheet2)
$ws9.Name = "t9_ML_2012"

for ($c = 1; $c -le $headers.Length; $c++) {
    $ws9.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$t9rows = @(
    @{ A="ERIPRET";          B="Paul";                           D="paul.eripret@etu.unilasalle.fr";                              E="13 min 37 s"; F=13;  G="2,00"; H="2,00"; I="0,00"; J="-";    K="3,00"; L="3,00"; M="-";    N="3,00" },
    @{ A="BOUTILLIER";       B="Hugo";                           D="hugo.boutillier@etu.unilasalle.fr";                           E="9 min 23 s";  F=9;   G="0,00"; H="0,00"; I="0,00"; J="0,00"; K="3,00"; L="3,00"; M="0,00"; N="3,00" },
    @{ A="LEKANGA MBOMA";    B="Amassa Roland Nathanael";        D="amassarolandnathanael.lekangamboma@etu.unilasalle.fr";       E="34 min 42 s"; F=8.5; G="0,00"; H="0,00"; I="2,50"; J="-";    K="0,00"; L="3,00"; M="-";    N="3,00" },
    @{ A="PECOURT";          B="Bertille";                       D="bertille.pecourt@etu.unilasalle.fr";                          E="4 min 44 s";  F=8;   G="0,00"; H="0,00"; I="0,00"; J="0,00"; K="0,00"; L="3,00"; M="2,00"; N="3,00" },
    @{ A="RAVELOJAONA";      B="Arthur";                         D="arthur.ravelojaona@etu.unilasalle.fr";                        E="22 min 28 s"; F=7;   G="0,00"; H="2,00"; I="0,00"; J="0,00"; K="0,00"; L="3,00"; M="2,00"; N="0,00" },
    @{ A="CHIEN-CHOW-CHINE"; B="Jules";                          D="jules.chien-chow-chine@etu.unilasalle.fr";                     E="3 min 4 s";   F=0;   G="0,00"; H="0,00"; I="0,00"; J="0,00"; K="0,00"; L="0,00"; M="0,00"; N="0,00" }
)

$r = 2
foreach ($row in $t9rows) {
    $ws9.Cells.Item($r, 1).Value = $row.A
    $ws9.Cells.Item($r, 2).Value = $row.B
    $ws9.Cells.Item($r, 3).Formula = "=A$r&B$r"
    $ws9.Cells.Item($r, 4).Value = $row.D
    $ws9.Cells.Item($r, 5).Value = $row.E
    $ws9.Cells.Item($r, 6).Value = $row.F
    $ws9.Cells.Item($r, 7).Value = $row.G
    $ws9.Cells.Item($r, 8).Value = $row.H
    $ws9.Cells.Item($r, 9).Value = $row.I
    $ws9.Cells.Item($r, 10).Value = $row.J
    $ws9.Cells.Item($r, 11).Value = $row.K
    $ws9.Cells.Item($r, 12).Value = $row.L
    $ws9.Cells.Item($r, 13).Value = $row.M
    $ws9.Cells.Item($r, 14).Value = $row.N
    $r++
}

[void]$ws9.Range("H22").Select()

# t9_ML_2012 is now the right-most / newest sheet and becomes the active tab,
# matching the source workbook's activeTab move.
$ws9.Activate()
